$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 5 (Stanbic IBTC Holdings) is dropped; the remaining companies'
# capital-structure metrics are refreshed with updated source data, and the
# company previously in row 3 (Abbey Mortgage Bank) is replaced by Stanbic
# IBTC Holdings with new figures. Delete the trailing row first so the sheet
# dimension shrinks from A1:AQ5 to A1:AQ4, then rewrite the per-row values.
$ws.Rows(5).Delete()

# Row 2 (Nigeria / company_name "3" -> "2")
# B2 must stay a text cell (it looks numeric), so force text format before
# assigning, matching how Excel would store a user-entered text "2".
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("D2").Value = 0.12595
$ws.Range("E2").Value = 0.2121
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 218.6
$ws.Range("L2").Value = 0.4239556262363756
$ws.Range("M2").Value = 48.7
$ws.Range("N2").Value = 0.03754529334669648
$ws.Range("O2").Value = 0.2227813357731016
$ws.Range("P2").Value = 48.7
$ws.Range("Q2").Value = 0.03754529334669648
$ws.Range("R2").Value = 0.2227813357731016
$ws.Range("U2").Value = 2434.33
$ws.Range("V2").Value = 1.876748130444838
$ws.Range("W2").Value = 0.171991967364791
$ws.Range("X2").Value = 0.03084704100146465
$ws.Range("Y2").Value = 0.1411449263633263
$ws.Range("Z2").Value = 2.992744790759766
$ws.Range("AB2").Value = 0.03701957132948389
$ws.Range("AC2").Value = -0.03701957132948389
$ws.Range("AD2").Value = 585.6799999999999
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 585.6799999999999
$ws.Range("AG2").Value = -1848.65
$ws.Range("AH2").Value = 0.3110719255568892
$ws.Range("AI2").Value = 0.3761399542733835
$ws.Range("AJ2").Value = 3.351736016680265
$ws.Range("AK2").Value = 2.107324023938444
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3 (company changes from Abbey Mortgage Bank Plc to Stanbic IBTC Holdings PLC)
$ws.Range("B3").Value = "Stanbic IBTC Holdings PLC (NGSE:STANBIC)"
$ws.Range("D3").Value = 0.171
$ws.Range("E3").Value = 0.337
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 217.6
$ws.Range("L3").Value = 0.4243369734789392
$ws.Range("M3").Value = 48.7
$ws.Range("N3").Value = 0.03798159413508033
$ws.Range("O3").Value = 0.2238051470588235
$ws.Range("P3").Value = 48.7
$ws.Range("Q3").Value = 0.03798159413508033
$ws.Range("R3").Value = 0.2238051470588235
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 2431.1
$ws.Range("V3").Value = 1.896038059585088
$ws.Range("W3").Value = 0.2745394902851375
$ws.Range("X3").Value = 0.03102264652923829
$ws.Range("Y3").Value = 0.2435168437558992
$ws.Range("Z3").Value = 3.373684210526315
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03614563170671756
$ws.Range("AC3").Value = -0.03614563170671756
$ws.Range("AD3").Value = 579.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 579.4
$ws.Range("AG3").Value = -1851.7
$ws.Range("AH3").Value = 0.3112376450365277
$ws.Range("AI3").Value = 0.3774346948081558
$ws.Range("AJ3").Value = 3.251448639157156
$ws.Range("AK3").Value = 2.066629464285715
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# Row 4 (Infinity Trust Mortgage Bank Plc, figures refreshed)
$ws.Range("D4").Value = 0.0809
$ws.Range("E4").Value = 0.0872
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3546099290780142
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 3.23
$ws.Range("V4").Value = 0.2167785234899329
$ws.Range("W4").Value = 0.06944444444444445
$ws.Range("X4").Value = 0.03067143547369102
$ws.Range("Y4").Value = 0.03877300897075343
$ws.Range("Z4").Value = 0.1389847215377033
$ws.Range("AB4").Value = 0.03789351095225023
$ws.Range("AC4").Value = -0.03789351095225023
$ws.Range("AD4").Value = 6.28
$ws.Range("AF4").Value = 6.28
$ws.Range("AG4").Value = 3.05
$ws.Range("AH4").Value = 0.2965061378659112
$ws.Range("AI4").Value = 0.2857142857142857
$ws.Range("AJ4").Value = 0.16991643454039
$ws.Range("AK4").Value = 0.1626666666666667
$ws.Range("T4").ClearContents()
